# Remove the header row (row 1: Salesperson/Date/Revenue/Cost) from the
# "Demo Sales Data" sheet, shifting all data rows up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(1).Delete()

# Match the resulting selection shown in the edited workbook.
$ws.Range("B8").Select() | Out-Null
